$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 717.8461
$ws.Range("I92").Value = 361
$ws.Range("J92").Value = 5000
$ws.Range("K92").Value = 361
$ws.Range("L92").Value = 5000
$ws.Range("M92").Value = 887
$ws.Range("N92").Value = -7496
$ws.Range("H113").Value = 10950.5
$ws.Range("I113").Value = 4166.6665
$ws.Range("K113").Value = 4166.6665
$ws.Range("M113").Value = -912.6665000000003
$ws.Range("H114").Value = 37865
$ws.Range("J114").Value = 37865
$ws.Range("L114").Value = 37865
$ws.Range("N114").Value = -46543
$ws.Range("H116").Value = 633637
$ws.Range("I116").Value = 1251810.8
$ws.Range("J116").Value = 15463.25
$ws.Range("K116").Value = 1251810.8
$ws.Range("L116").Value = 15463.25
$ws.Range("M116").Value = -1248368.8
$ws.Range("N116").Value = -22347.25
$ws.Range("H138").Value = 2790.45
$ws.Range("J138").Value = 2954.6965
$ws.Range("L138").Value = 8864.0895
$ws.Range("N138").Value = -19144.0895

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6085.7
$ws.Range("I32").Value = 3375.2222
$ws.Range("K32").Value = 3375.2222
$ws.Range("M32").Value = -3088.2222
$ws.Range("H74").Value = 6396.381
$ws.Range("I74").Value = 6450.8237
$ws.Range("J74").Value = 6165
$ws.Range("K74").Value = 6450.8237
$ws.Range("L74").Value = 6165
$ws.Range("M74").Value = -5576.8237
$ws.Range("N74").Value = -7913
$ws.Range("H77").Value = 6396.381
$ws.Range("I77").Value = 6450.8237
$ws.Range("J77").Value = 6165
$ws.Range("K77").Value = 32254.1185
$ws.Range("L77").Value = 30825
$ws.Range("M77").Value = -27886.1185
$ws.Range("N77").Value = -39561
$ws.Range("H109").Value = 26050
$ws.Range("J109").Value = 26050
$ws.Range("L109").Value = 26050
$ws.Range("N109").Value = -28824
$ws.Range("H122").Value = 2686.7222
$ws.Range("I122").Value = 1381.5714
$ws.Range("K122").Value = 4144.7142
$ws.Range("M122").Value = -1694.7142
$ws.Range("H132").Value = 2347.0386
$ws.Range("I132").Value = 1639.2972
$ws.Range("J132").Value = 4092.8
$ws.Range("K132").Value = 4917.8916
$ws.Range("L132").Value = 12278.4
$ws.Range("M132").Value = -2387.8916
$ws.Range("N132").Value = -17338.4
$ws.Range("H137").Value = 41791.6
$ws.Range("J137").Value = 41791.6
$ws.Range("L137").Value = 41791.6
$ws.Range("N137").Value = -51991.6

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 467
$ws.Range("I22").Value = 467
$ws.Range("K22").Value = 467
$ws.Range("M22").Value = -294
$ws.Range("H94").Value = 855.9355
$ws.Range("I94").Value = 923.8182
$ws.Range("J94").Value = 690
$ws.Range("K94").Value = 923.8182
$ws.Range("L94").Value = 690
$ws.Range("M94").Value = -472.8182
$ws.Range("N94").Value = -1592
$ws.Range("H125").Value = 41776.668
$ws.Range("J125").Value = 41776.668
$ws.Range("L125").Value = 41776.668
$ws.Range("N125").Value = -51616.668
$ws.Range("H137").Value = 40892
$ws.Range("J137").Value = 40892
$ws.Range("L137").Value = 40892
$ws.Range("N137").Value = -51092

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 603.6429000000001
$ws.Range("I22").Value = 211.22223
$ws.Range("J22").Value = 1310
$ws.Range("K22").Value = 211.22223
$ws.Range("L22").Value = 1310
$ws.Range("M22").Value = 138.77777
$ws.Range("N22").Value = -2010
$ws.Range("H31").Value = 19234176
$ws.Range("I31").Value = 2076.611
$ws.Range("J31").Value = 62506400
$ws.Range("K31").Value = 2076.611
$ws.Range("L31").Value = 62506400
$ws.Range("M31").Value = -1781.611
$ws.Range("N31").Value = -62506990
$ws.Range("H34").Value = 19234176
$ws.Range("I34").Value = 2076.611
$ws.Range("J34").Value = 62506400
$ws.Range("K34").Value = 2076.611
$ws.Range("L34").Value = 62506400
$ws.Range("M34").Value = -1874.611
$ws.Range("N34").Value = -62506804
$ws.Range("H55").Value = 25599
$ws.Range("J55").Value = 25599
$ws.Range("L55").Value = 25599
$ws.Range("N55").Value = -26229
$ws.Range("H58").Value = 1979.2858
$ws.Range("I58").Value = 1680.7192
$ws.Range("J58").Value = 4815.6665
$ws.Range("K58").Value = 1680.7192
$ws.Range("L58").Value = 4815.6665
$ws.Range("M58").Value = -1477.7192
$ws.Range("N58").Value = -5221.6665
$ws.Range("H81").Value = 34900
$ws.Range("J81").Value = 34900
$ws.Range("L81").Value = 34900
$ws.Range("N81").Value = -36896
$ws.Range("H84").Value = 34900
$ws.Range("J84").Value = 34900
$ws.Range("L84").Value = 104700
$ws.Range("N84").Value = -114684
$ws.Range("H132").Value = 2283.3125
$ws.Range("I132").Value = 1008.10345
$ws.Range("J132").Value = 4229.684
$ws.Range("K132").Value = 3024.31035
$ws.Range("L132").Value = 12689.052
$ws.Range("M132").Value = -494.3103499999997
$ws.Range("N132").Value = -17749.052
$ws.Range("H136").Value = 1979.2858
$ws.Range("I136").Value = 1680.7192
$ws.Range("J136").Value = 4815.6665
$ws.Range("K136").Value = 5042.1576
$ws.Range("L136").Value = 14446.9995
$ws.Range("M136").Value = -2492.1576
$ws.Range("N136").Value = -19546.9995

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1523.9615
$ws.Range("I5").Value = 296.75
$ws.Range("J5").Value = 5614.6665
$ws.Range("K5").Value = 890.25
$ws.Range("L5").Value = 16843.9995
$ws.Range("M5").Value = -778.25
$ws.Range("N5").Value = -17067.9995
$ws.Range("H107").Value = 51103.95
$ws.Range("I107").Value = 377.42856
$ws.Range("J107").Value = 169465.83
$ws.Range("K107").Value = 1132.28568
$ws.Range("L107").Value = 508397.49
$ws.Range("M107").Value = 787.71432
$ws.Range("N107").Value = -512237.49
$ws.Range("H132").Value = 2181.9
$ws.Range("I132").Value = 843.4286
$ws.Range("J132").Value = 5305
$ws.Range("K132").Value = 7590.8574
$ws.Range("L132").Value = 47745
$ws.Range("M132").Value = -5060.8574
$ws.Range("N132").Value = -52805
$ws.Range("H135").Value = 1523.9615
$ws.Range("I135").Value = 296.75
$ws.Range("J135").Value = 5614.6665
$ws.Range("K135").Value = 2670.75
$ws.Range("L135").Value = 50531.9985
$ws.Range("M135").Value = -135.75
$ws.Range("N135").Value = -55601.9985

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 71.25
$ws.Range("I2").Value = 62.5
$ws.Range("K2").Value = 62.5
$ws.Range("M2").Value = 50.5
$ws.Range("H102").Value = 2282.3845
$ws.Range("I102").Value = 1681.5264
$ws.Range("J102").Value = 3913.2856
$ws.Range("K102").Value = 1681.5264
$ws.Range("L102").Value = 3913.2856
$ws.Range("M102").Value = -59.52639999999997
$ws.Range("N102").Value = -7157.2856
$ws.Range("H122").Value = 4321.231
$ws.Range("I122").Value = 2616.8
$ws.Range("K122").Value = 7850.400000000001
$ws.Range("M122").Value = -5400.400000000001
$ws.Range("H128").Value = 42351.43
$ws.Range("J128").Value = 42351.43
$ws.Range("L128").Value = 42351.43
$ws.Range("N128").Value = -52311.43
$ws.Range("H137").Value = 74212.25
$ws.Range("J137").Value = 74212.25
$ws.Range("L137").Value = 74212.25
$ws.Range("N137").Value = -84412.25

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 44600
$ws.Range("J75").Value = 44600
$ws.Range("L75").Value = 44600
$ws.Range("N75").Value = -46472
$ws.Range("H78").Value = 44600
$ws.Range("J78").Value = 44600
$ws.Range("L78").Value = 133800
$ws.Range("N78").Value = -143160
$ws.Range("I122").Value = 2847.9
$ws.Range("J122").Value = 8285
$ws.Range("K122").Value = 8543.700000000001
$ws.Range("L122").Value = 24855
$ws.Range("M122").Value = -6093.700000000001
$ws.Range("N122").Value = -29755

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 44619.8
$ws.Range("J80").Value = 44619.8
$ws.Range("L80").Value = 44619.8
$ws.Range("N80").Value = -46615.8
$ws.Range("H83").Value = 44619.8
$ws.Range("J83").Value = 44619.8
$ws.Range("L83").Value = 133859.4
$ws.Range("N83").Value = -143843.4
$ws.Range("H108").Value = 39550
$ws.Range("J108").Value = 39550
$ws.Range("L108").Value = 39550
$ws.Range("N108").Value = -47230
$ws.Range("H122").Value = 3211.7058
$ws.Range("I122").Value = 2049.9167
$ws.Range("J122").Value = 6000
$ws.Range("K122").Value = 6149.750100000001
$ws.Range("L122").Value = 18000
$ws.Range("M122").Value = -3699.750100000001
$ws.Range("N122").Value = -22900
$ws.Range("H132").Value = 7408989
$ws.Range("I132").Value = 570.4231
$ws.Range("J132").Value = 17546824
$ws.Range("K132").Value = 1711.2693
$ws.Range("L132").Value = 52640472
$ws.Range("M132").Value = 818.7307000000001
$ws.Range("N132").Value = -52645532
$ws.Range("H136").Value = 1463.36
$ws.Range("I136").Value = 660.58826
$ws.Range("J136").Value = 3169.25
$ws.Range("K136").Value = 1981.76478
$ws.Range("L136").Value = 9507.75
$ws.Range("M136").Value = 568.23522
$ws.Range("N136").Value = -14607.75
